# QRDecompTest.xlsx edit:
#  - add a new blank worksheet ("Sheet1") at the end of the tab strip
#  - re-select the "6x6" sheet (it stays the active tab) and move the
#    cell cursor to K3 (previously a full A1:XFD8 row-header selection)
#  - RAND()-driven "rand" sheet recalculates on its own (volatile formulas)

$wb = $excel.ActiveWorkbook

# Append a brand-new worksheet after the last existing tab so it lands
# at the end of the sheet list (Sheets.Add with no args would insert
# before the active sheet instead).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Jump back to the "6x6" sheet (the tab that was active before/after the
# edit) and move the selection to K3.
$sixBySix = $wb.Worksheets.Item("6x6")
$sixBySix.Activate() | Out-Null
$sixBySix.Range("K3").Select() | Out-Null
